# Generate Report for Handoff
# Updates the status of "b.md" (row 3) on the Overview sheet and the
# per-language detail sheets (zh-cn, de-de) to reflect that a new handoff
# package was generated for it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md's zh-cn / de-de status + the "Latest HO Xliff
# Generate Date" column move from the old "Handed back" state to
# "Ready for handoff" with the new generation timestamp.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 16:45:22"

# ---------------------------------------------------------------------
# zh-cn detail sheet: row 3 is b.md's record.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 16:45:14"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/47b776dd3870d9a7c2310472221bee63b509bf33/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de detail sheet: row 3 is b.md's record.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 16:45:22"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcb2561f9fb3043c3c67aa4a6fb3beced90cae95/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/47b776dd3870d9a7c2310472221bee63b509bf33/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 40
